$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3296.4443
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 4502
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 4502
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -5132
$ws.Range("H79").Value = 3296.4443
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 4502
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 4502
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -6686
$ws.Range("H125").Value = 2232.4
$ws.Range("I125").Value = 2232.4
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 20091.6
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -17631.6
$ws.Range("N125").ClearContents()
$ws.Range("H129").Value = 909.3261
$ws.Range("J129").Value = 962.6667
$ws.Range("L129").Value = 2888.0001
$ws.Range("N129").Value = -12888.0001
$ws.Range("H132").Value = 95074.72
$ws.Range("I132").Value = 115815.23
$ws.Range("K132").Value = 347445.69
$ws.Range("M132").Value = -344915.69

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 745.1667
$ws.Range("I2").Value = 691.6
$ws.Range("K2").Value = 691.6
$ws.Range("M2").Value = -578.6
$ws.Range("H32").Value = 17129.428
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 17129.428
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 17129.428
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -17703.428
$ws.Range("H35").Value = 24444.334
$ws.Range("I35").Value = 4000
$ws.Range("J35").Value = 34666.5
$ws.Range("K35").Value = 4000
$ws.Range("L35").Value = 34666.5
$ws.Range("M35").Value = -3594
$ws.Range("N35").Value = -35478.5
$ws.Range("H61").Value = 1961.3889
$ws.Range("I61").Value = 1129.5454
$ws.Range("J61").Value = 3268.5715
$ws.Range("K61").Value = 1129.5454
$ws.Range("L61").Value = 3268.5715
$ws.Range("M61").Value = -917.5454
$ws.Range("N61").Value = -3692.5715
$ws.Range("H63").Value = 9896714
$ws.Range("I63").Value = 12594400
$ws.Range("K63").Value = 12594400
$ws.Range("M63").Value = -12593714
$ws.Range("H66").Value = 9896714
$ws.Range("I66").Value = 12594400
$ws.Range("K66").Value = 62972000
$ws.Range("M66").Value = -62968568
$ws.Range("H116").Value = 745.1667
$ws.Range("I116").Value = 691.6
$ws.Range("K116").Value = 691.6
$ws.Range("M116").Value = 1602.4
$ws.Range("H136").Value = 1961.3889
$ws.Range("I136").Value = 1129.5454
$ws.Range("J136").Value = 3268.5715
$ws.Range("K136").Value = 3388.6362
$ws.Range("L136").Value = 9805.7145
$ws.Range("M136").Value = -838.6361999999999
$ws.Range("N136").Value = -14905.7145

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 745.1667
$ws.Range("I3").Value = 691.6
$ws.Range("K3").Value = 691.6
$ws.Range("M3").Value = -577.6
$ws.Range("H8").Value = 1280.2307
$ws.Range("I8").Value = 655.4
$ws.Range("J8").Value = 3363
$ws.Range("K8").Value = 655.4
$ws.Range("L8").Value = 3363
$ws.Range("M8").Value = -515.4
$ws.Range("N8").Value = -3643
$ws.Range("H105").Value = 2492.2307
$ws.Range("I105").Value = 2445.3635
$ws.Range("K105").Value = 2445.3635
$ws.Range("M105").Value = -698.3634999999999

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4179.56
$ws.Range("I31").Value = 2176.625
$ws.Range("J31").Value = 5122.1177
$ws.Range("K31").Value = 2176.625
$ws.Range("L31").Value = 5122.1177
$ws.Range("M31").Value = -1881.625
$ws.Range("N31").Value = -5712.1177
$ws.Range("H34").Value = 4179.56
$ws.Range("I34").Value = 2176.625
$ws.Range("J34").Value = 5122.1177
$ws.Range("K34").Value = 2176.625
$ws.Range("L34").Value = 5122.1177
$ws.Range("M34").Value = -1974.625
$ws.Range("N34").Value = -5526.1177
$ws.Range("H58").Value = 2147.88
$ws.Range("I58").Value = 1885.5918
$ws.Range("K58").Value = 1885.5918
$ws.Range("M58").Value = -1682.5918
$ws.Range("H135").Value = 49045
$ws.Range("J135").Value = 49045
$ws.Range("L135").Value = 49045
$ws.Range("N135").Value = -59185
$ws.Range("H136").Value = 2147.88
$ws.Range("I136").Value = 1885.5918
$ws.Range("K136").Value = 5656.7754
$ws.Range("M136").Value = -3106.7754
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980
$ws.Range("H139").Value = 48035
$ws.Range("J139").Value = 48035
$ws.Range("L139").Value = 48035
$ws.Range("N139").Value = -58315

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H92").Value = 33335780
$ws.Range("J92").Value = 62504156
$ws.Range("L92").Value = 187512468
$ws.Range("N92").Value = -187514964
$ws.Range("H113").Value = 931.0417
$ws.Range("I113").Value = 742.25
$ws.Range("J113").Value = 1875
$ws.Range("K113").Value = 2226.75
$ws.Range("L113").Value = 5625
$ws.Range("M113").Value = -56.75
$ws.Range("N113").Value = -9965
$ws.Range("H131").Value = 9616390
$ws.Range("I131").Value = 35715480
$ws.Range("J131").Value = 935.7368
$ws.Range("K131").Value = 107146440
$ws.Range("L131").Value = 2807.2104
$ws.Range("M131").Value = -107141400
$ws.Range("N131").Value = -12887.2104

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 3042.8572
$ws.Range("I17").Value = 2800
$ws.Range("J17").Value = 4500
$ws.Range("K17").Value = 2800
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = -2632
$ws.Range("N17").Value = -4836
$ws.Range("H70").Value = 6273.593
$ws.Range("I70").Value = 5564.35
$ws.Range("K70").Value = 5564.35
$ws.Range("M70").Value = -5294.35
$ws.Range("H73").Value = 6273.593
$ws.Range("I73").Value = 5564.35
$ws.Range("K73").Value = 5564.35
$ws.Range("M73").Value = -4628.35
$ws.Range("H80").Value = 41668900
$ws.Range("I80").Value = 62501850
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 62501850
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -62500852
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 41668900
$ws.Range("I83").Value = 62501850
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 312509250
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -312504258
$ws.Range("N83").Value = -24984

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 4671967.5
$ws.Range("I11").Value = 2900
$ws.Range("J11").Value = 7006501
$ws.Range("K11").Value = 2900
$ws.Range("L11").Value = 7006501
$ws.Range("M11").Value = -2760
$ws.Range("N11").Value = -7006781
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H136").Value = 2683.6829
$ws.Range("I136").Value = 1416.7407
$ws.Range("J136").Value = 5127.0713
$ws.Range("K136").Value = 4250.2221
$ws.Range("L136").Value = 15381.2139
$ws.Range("M136").Value = -1700.2221
$ws.Range("N136").Value = -20481.2139

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 97231.39999999999
$ws.Range("J46").Value = 97231.39999999999
$ws.Range("L46").Value = 97231.39999999999
$ws.Range("N46").Value = -97693.39999999999
$ws.Range("H134").Value = 97231.39999999999
$ws.Range("J134").Value = 97231.39999999999
$ws.Range("L134").Value = 291694.2
$ws.Range("N134").Value = -296764.2
$ws.Range("H136").Value = 4551.16
$ws.Range("I136").Value = 1591.1
$ws.Range("J136").Value = 6524.533
$ws.Range("K136").Value = 4773.299999999999
$ws.Range("L136").Value = 19573.599
$ws.Range("M136").Value = -2223.299999999999
$ws.Range("N136").Value = -24673.599
